# Use solver to find minimum number of products needed to reach 2400 calories
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calories-Solver")
$ws.Activate()

# --- New "How many?" decision-variable column -----------------------------
$ws.Range("C1").Value = "How many?"

# Quantities chosen by Solver: 1x Licorice Rope, 2x Nachos, 2x Popcorn = 2400 cal
$qty = @(0, 0, 0, 0, 0, 0, 0, 0, 1, 2, 0, 2, 0, 0)
for ($i = 0; $i -lt $qty.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 3).Value = $qty[$i]
}

# --- Summary rows -----------------------------------------------------------
$ws.Range("B16").Value = "Total Items"
$ws.Range("C16").Formula = "=SUM(C2:C15)"

$ws.Range("B17").Value = "Total Calories"
$ws.Range("C17").Formula = "=SUMPRODUCT(B2:B15,C2:C15)"

# Match the saved selection / page orientation left behind by Solver's run
$ws.Range("D17").Select()
$ws.PageSetup.Orientation = 1

# --- Persist the Solver model parameters as hidden workbook-scoped (sheet) names
function Add-SolverName($name, $refersTo) {
    $n = $ws.Names.Add($name, $refersTo)
    $n.Visible = $false
}

Add-SolverName 'solver_adj'  '=''Calories-Solver''!$C$2:$C$15'
Add-SolverName 'solver_cvg'  '=0.0001'
Add-SolverName 'solver_drv'  '=1'
Add-SolverName 'solver_eng'  '=2'
Add-SolverName 'solver_itr'  '=2147483647'
Add-SolverName 'solver_lhs1' '=''Calories-Solver''!$C$17'
Add-SolverName 'solver_lhs2' '=''Calories-Solver''!$C$2:$C$15'
Add-SolverName 'solver_lin'  '=1'
Add-SolverName 'solver_mip'  '=2147483647'
Add-SolverName 'solver_mni'  '=30'
Add-SolverName 'solver_mrt'  '=0.075'
Add-SolverName 'solver_msl'  '=2'
Add-SolverName 'solver_neg'  '=1'
Add-SolverName 'solver_nod'  '=2147483647'
Add-SolverName 'solver_num'  '=2'
Add-SolverName 'solver_opt'  '=''Calories-Solver''!$C$16'
Add-SolverName 'solver_pre'  '=0.000001'
Add-SolverName 'solver_rbv'  '=1'
Add-SolverName 'solver_rel1' '=2'
Add-SolverName 'solver_rel2' '=4'
Add-SolverName 'solver_rhs1' '=2400'
Add-SolverName 'solver_rhs2' '=integer'
Add-SolverName 'solver_rlx'  '=2'
Add-SolverName 'solver_rsd'  '=0'
Add-SolverName 'solver_scl'  '=2'
Add-SolverName 'solver_sho'  '=2'
Add-SolverName 'solver_ssz'  '=100'
Add-SolverName 'solver_tim'  '=2147483647'
Add-SolverName 'solver_tol'  '=0.01'
Add-SolverName 'solver_typ'  '=2'
Add-SolverName 'solver_val'  '=0'
Add-SolverName 'solver_ver'  '=2'

$wb.Calculate()
